$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 44907
$ws.Range("F2").Value = 0.0135406646415844
$ws.Range("G2").Value = 1.00619485639197
$ws.Range("I2").Value = 0.00196336990767471
$ws.Range("J2").Value = 0.00861649090764852
$ws.Range("K2").Value = 0.000414262513320785

# Row 3
$ws.Range("E3").Value = 44907
$ws.Range("F3").Value = 0.103935692392368
$ws.Range("G3").Value = 1.39045246368396
$ws.Range("I3").Value = 0.0457945946206097
$ws.Range("J3").Value = 0.137297394402522
$ws.Range("K3").Value = 0.0137871562209999

# Row 4
$ws.Range("E4").Value = 44907
$ws.Range("F4").Value = 0.112483257021618
$ws.Range("G4").Value = 2.40966284884659
$ws.Range("H4").Value = -0.0119310649580203
$ws.Range("I4").Value = 0.0392412304075698
$ws.Range("J4").Value = 0.144118176905062
$ws.Range("K4").Value = 0.00598365806940942

# Row 5
$ws.Range("E5").Value = 44907
$ws.Range("F5").Value = 0.290083895429395
$ws.Range("G5").Value = 14.2222222222222
$ws.Range("I5").Value = 0.186934849799873
$ws.Range("J5").Value = 0.392035473943881
$ws.Range("K5").Value = 0.0722292075782106

# Row 6
$ws.Range("D6").Value = 43878
$ws.Range("E6").Value = 44907
$ws.Range("F6").Value = 0.135210958198813
$ws.Range("G6").Value = 1.74235782893734
$ws.Range("I6").Value = 0.0290100042592723
$ws.Range("J6").Value = 0.166449080132706
$ws.Range("K6").Value = 0.0013281937486539
